$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ 2 = 10.25756710338925; 3 = 5.342309548077196; 4 = 9.297245342656401; 5 = 13.7335740058896; 6 = 32.52346818206429; 8 = 7.344005520526261; 9 = 22.24974157687664; 10 = 9.921764148115141; 11 = 10.47487306526676; 14 = 18.85652582045357; 15 = 24.48923716574428 }
    3 = @{ 2 = 9.984085959012209; 3 = 5.147060776718823; 4 = 9.234468675093204; 5 = 13.66831861141322; 6 = 32.56530234486046; 8 = 7.344005520526261; 9 = 22.33373777365707; 10 = 9.927571469992708; 11 = 10.29322046947316; 14 = 18.91345717709838; 15 = 24.55969453908657 }
    4 = @{ 2 = 9.813959927715436; 3 = 5.024042825955766; 4 = 9.197417890228213; 5 = 13.63099248755583; 6 = 32.59838166921292; 8 = 7.344005520526261; 9 = 22.38922204329311; 10 = 9.932686044028589; 11 = 10.18169473288713; 14 = 18.95006846841527; 15 = 24.60776970458099 }
    5 = @{ 2 = 9.7441834292288; 3 = 4.973201039009398; 4 = 9.182708343453506; 5 = 13.61648340286937; 6 = 32.61371867737729; 8 = 7.344005520526261; 9 = 22.41281514391227; 10 = 9.935160187146936; 11 = 10.13630764894357; 14 = 18.96540528594443; 15 = 24.62856941878228 }
    6 = @{ 2 = 9.73257320973204; 3 = 4.964718264781549; 4 = 9.180289699008878; 5 = 13.61411689918649; 6 = 32.61637745453368; 8 = 7.344005520526261; 9 = 22.416792106563; 10 = 9.93559457725898; 11 = 10.12877652901302; 14 = 18.96797719933506; 15 = 24.63209614194815 }
    7 = @{ 2 = 9.813020569770766; 3 = 5.023359926888387; 4 = 9.197217920337545; 5 = 13.63079395650545; 6 = 32.59858099417733; 8 = 7.344005520526261; 9 = 22.38953624933868; 10 = 9.932717831942437; 11 = 10.18108230425521; 14 = 18.95027361439146; 15 = 24.60804532535969 }
    8 = @{ 2 = 10.16379096040415; 3 = 5.275682272926254; 4 = 9.275297890751492; 5 = 13.71051277203896; 6 = 32.53635710820753; 8 = 7.344005520526261; 9 = 22.27789161465178; 10 = 9.923445341818041; 11 = 10.41227184375978; 14 = 18.87581294114668; 15 = 24.51253063482077 }
    9 = @{ 2 = 10.82976924833978; 3 = 5.742607939962267; 4 = 9.439662774317867; 5 = 13.8880027978997; 6 = 32.47306548654721; 8 = 7.344005520526261; 9 = 22.09000344145986; 10 = 9.9175287845464; 11 = 10.86320818063559; 14 = 18.74287334184325; 15 = 24.36350302943202 }
    10 = @{ 2 = 11.30033681783906; 3 = 6.065149207830128; 4 = 9.566432911740256; 5 = 14.03046189438895; 6 = 32.46242533222458; 8 = 7.344005520526261; 9 = 21.9709167265886; 10 = 9.92062454093602; 11 = 11.18978166969475; 14 = 18.65309574235725; 15 = 24.27745002746433 }
    11 = @{ 2 = 11.50934082294159; 3 = 6.206822168571673; 4 = 9.625220308046025; 5 = 14.09769235535218; 6 = 32.46536903938014; 8 = 7.344005520526261; 9 = 21.92086331595034; 10 = 9.923639104852251; 11 = 11.33667009771465; 14 = 18.61395060943256; 15 = 24.24341141619842 }
    12 = @{ 2 = 11.58768337372814; 3 = 6.259701194822862; 4 = 9.647626139655562; 5 = 14.12348195915714; 6 = 32.46760147230453; 8 = 7.344005520526261; 9 = 21.90250244880759; 10 = 9.925010559799112; 11 = 11.39200130176243; 14 = 18.59936988239864; 15 = 24.23125774153415 }
    13 = @{ 2 = 11.57084773038122; 3 = 6.248347595968555; 4 = 9.642794490392625; 5 = 14.11791328031119; 6 = 32.46707100132702; 8 = 7.344005520526261; 9 = 21.9064303902375; 10 = 9.924704985830539; 11 = 11.38009853395669; 14 = 18.602499326349; 15 = 24.23384249499174 }
    14 = @{ 2 = 11.51580260617686; 3 = 6.211188204740107; 4 = 9.627060852604943; 5 = 14.09980755430638; 6 = 32.46553031434635; 8 = 7.344005520526261; 9 = 21.91934086130616; 10 = 9.923747334250105; 11 = 11.34122837450082; 14 = 18.61274618658626; 15 = 24.24239676681063 }
    15 = @{ 2 = 11.48197922368368; 3 = 6.188325621815814; 4 = 9.617441844246123; 5 = 14.08875982641322; 6 = 32.46473209726106; 8 = 7.344005520526261; 9 = 21.92732618506323; 10 = 9.923190651327063; 11 = 11.31737972285759; 14 = 18.6190542574319; 15 = 24.24773239984489 }
    16 = @{ 2 = 11.28656952229423; 3 = 6.055785080066321; 4 = 9.562612120854869; 5 = 14.02611551331224; 6 = 32.46238946777735; 8 = 7.344005520526261; 9 = 21.97427080750683; 10 = 9.920459765479116; 11 = 11.18014403618528; 14 = 18.65568798740462; 15 = 24.27977744166421 }
    17 = @{ 2 = 11.16534047827763; 3 = 5.973149328450774; 4 = 9.529250730623426; 5 = 13.9882938794545; 6 = 32.46294510629769; 8 = 7.344005520526261; 9 = 22.00412558333349; 10 = 9.919195137530565; 11 = 11.09548854917291; 14 = 18.67859499361972; 15 = 24.30074531614047 }
    18 = @{ 2 = 11.09514056666594; 3 = 5.925146035354492; 4 = 9.510168477661475; 5 = 13.96676924577267; 6 = 32.46399750485921; 8 = 7.344005520526261; 9 = 22.02168497717166; 10 = 9.918619020832962; 11 = 11.04664257124749; 14 = 18.69193011652229; 15 = 24.31328612328838 }
    19 = @{ 2 = 11.07129338897521; 3 = 5.908813038213057; 4 = 9.503726320297789; 5 = 13.95952130302031; 6 = 32.46447972588719; 8 = 7.344005520526261; 9 = 22.02769685089963; 10 = 9.918449968118097; 11 = 11.03007927034586; 14 = 18.69647260764177; 15 = 24.31761472510023 }
    20 = @{ 2 = 11.17829495107112; 3 = 5.981995368168603; 4 = 9.532791221736476; 5 = 13.99229644542147; 6 = 32.46281012368404; 8 = 7.344005520526261; 9 = 22.00090735808562; 10 = 9.919314111504749; 11 = 11.10451663217206; 14 = 18.67613999157157; 15 = 24.29846349372521 }
    21 = @{ 2 = 11.53199304486849; 3 = 6.222124015631668; 4 = 9.63167841975819; 5 = 14.10511680972772; 6 = 32.4659525337562; 8 = 7.344005520526261; 9 = 21.91553263481163; 10 = 9.92402238938052; 11 = 11.35265381509554; 14 = 18.60972985744817; 15 = 24.23986418286232 }
    22 = @{ 2 = 11.75844612987035; 3 = 6.374560937936613; 4 = 9.697140319404681; 5 = 14.18077144616164; 6 = 32.4745198705791; 8 = 7.344005520526261; 9 = 21.86319377632709; 10 = 9.928438995859421; 11 = 11.51310009458757; 14 = 18.56774106048747; 15 = 24.20585679921923 }
    23 = @{ 2 = 11.63803783299685; 3 = 6.293627087501326; 4 = 9.662131341560801; 5 = 14.14022349853714; 6 = 32.46935206470013; 8 = 7.344005520526261; 9 = 21.89081127323377; 10 = 9.9259595959218; 11 = 11.42764127074018; 14 = 18.59002223679715; 15 = 24.22361409963854 }
    24 = @{ 2 = 11.17243979710239; 3 = 5.977997611796996; 4 = 9.531190260269236; 5 = 13.99048620020815; 6 = 32.46286886610261; 8 = 7.344005520526261; 9 = 22.00236108586131; 10 = 9.919259853146849; 11 = 11.10043558134264; 14 = 18.67724938213366; 15 = 24.29949359158735 }
    25 = @{ 2 = 10.65252959041884; 3 = 5.619666703328003; 4 = 9.39407897981306; 5 = 13.83780537589692; 6 = 32.48389137838548; 8 = 7.344005520526261; 9 = 22.13750560257454; 10 = 9.917819324322547; 11 = 10.74182011742914; 14 = 18.77744522385475; 15 = 24.3997100659495 }
}

foreach ($row in $data.Keys) {
    foreach ($col in $data[$row].Keys) {
        $ws.Cells.Item($row, $col).Value = $data[$row][$col]
    }
}

Write-Host "Updated $($data.Count) rows"
